$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices and 1h volume % changes).
# Numeric-looking text values are written via a Text-number-format round trip
# so Excel keeps them as literal strings (matching the source inline strings)
# instead of silently parsing them into floating point numbers.

$ws.Range("D2").Value = "88.788.56"
$ws.Range("E2").Value = "  +9.09%  "
$ws.Range("D3").Value = "3.343.14"
$ws.Range("E3").Value = "  +5.26%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.393"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +34.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.29%  "
$ws.Range("D10").Value = "3.340.70"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.588"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000288"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.27%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.20%  "
$ws.Range("D15").Value = "3.952.50"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "88.516.87"
$ws.Range("E17").Value = "  +8.96%  "
$ws.Range("D18").Value = "3.337.55"
$ws.Range("E18").Value = "  +5.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +14.10%  "
$ws.Range("D26").Value = "3.508.65"
$ws.Range("E26").Value = "  +4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "78.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000130"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("B29").Value = "Cronos"
$ws.Range("C29").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.200"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +44.05%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "606.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("E33").Value = "  +6.95%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +23.15%  "
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +7.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.423"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "159.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "191.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +7.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.788"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.664"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.79%  "
